$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "27.728.91", "1.124")
# but must stay plain text, since "." is used both as thousands AND decimal
# separator in the source data (matching the original inlineStr content).
# Setting NumberFormat to text ("@") right before assigning the value prevents
# Excel from re-interpreting the string as a number. The style is then reset
# back to "Normal" immediately after so the cell keeps its original (default)
# formatting, matching cells that never had a custom style in the workbook.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.728.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.73%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.757.71"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.96%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.90%  "

$ws.Range("E6").Value = "  -0.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4434"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.32%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3734"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.41"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.74%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07542"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.16%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.124"
$ws.Range("D11").Style = "Normal"

$ws.Range("E12").Value = "  -0.24%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.203"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.19%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.425"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.757.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.99%  "

$ws.Range("E17").Value = "  -1.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.43"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +9.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06218"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.56%  "

$ws.Range("E20").Value = "  -0.08%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.32"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.41%  "

$ws.Range("E22").Value = "  -2.82%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5328"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "27.758.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.63%  "

$ws.Range("E25").Value = "  -0.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.321"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.08%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.11%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.361"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.956.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "128.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.59%  "

$ws.Range("E32").Value = "  -0.78%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09351"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.757"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.45%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.653"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -9.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.67"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02330"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2179"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.85%  "

$ws.Range("E39").Value = "  -2.58%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6484"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.86%  "

$ws.Range("E41").Value = "  -2.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.204"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.010"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.49%  "

$ws.Range("E44").Value = "  -4.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.08%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.60%  "

$ws.Range("E47").Value = "  -0.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.756"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.97%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "126.29"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.992"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06900"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.17%  "

Write-Output "Updated cryptos list"